# Generate Report for Archive
#
# The file "88f0e60b-118d-45c5-888e-a66eae44c829.md" (row 4 on every sheet)
# moves from "Ready for handoff" status to "In Translation" status, since a
# new handoff report has been generated / archived for it. Update the
# "Status" column on each sheet accordingly.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: File Name / zh-cn / de-de columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"

# --- zh-cn sheet: Status column (B) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B4").Value = "In Translation"

# --- de-de sheet: Status column (B) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B4").Value = "In Translation"
